# Updated the running times in the analysis
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Day 7's measured running times were updated (Initialization / Part 1 / Part 2)
$ws.Range("B8").Value = 2.0331999999999999
$ws.Range("C8").Value = 1.9512
$ws.Range("D8").Value = 0.0124

# Leave the selection where the user ended up after editing the row
$ws.Activate()
$ws.Range("D9").Select()

# Best-effort: make sure the chart that plots these columns picks up the
# refreshed numbers (no-op / silently ignored on hosts that already
# recompute chart caches automatically on cell change).
try {
    foreach ($co in $ws.ChartObjects()) {
        $co.Chart.Refresh()
    }
} catch {
}
